$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = 1.29
    "H2"  = 5.3
    "I2"  = 8.75
    "K2"  = 2.67
    "L2"  = 7
    "N2"  = 9.75
    "P2"  = 5
    "Q2"  = 1.44
    "R2"  = 2.6
    "S2"  = 1.24
    "T2"  = 3.65
    "X2"  = 7.4
    "Z2"  = 8.5
    "AC2" = 9.75
    "AD2" = 11
    "AE2" = 19
    "AF2" = 70
    "AG2" = 450
    "AH2" = 30
    "AJ2" = 27
    "AL2" = 90
    "AM2" = 65
    "AN2" = 3.35
    "AT2" = 3.65
    "AU2" = 8
    "AV2" = 60
    "AW2" = 9.75
    "AX2" = 45
    "AY2" = 37
    "AZ2" = 300
    "BA2" = 250
    "BB2" = 400
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
